$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Hunk 0 (row group)
$ws.Range("H86").Value = 1863.091
$ws.Range("J86").Value = 2668
$ws.Range("L86").Value = 2668
$ws.Range("N86").Value = -4914
# Hunk 1 (row group)
$ws.Range("H89").Value = 1863.091
$ws.Range("J89").Value = 2668
$ws.Range("L89").Value = 13340
$ws.Range("N89").Value = -24572
# Hunk 2 (row group)
$ws.Range("H139").Value = 39953.156
$ws.Range("J139").Value = 39953.156
$ws.Range("L139").Value = 39953.156
$ws.Range("N139").Value = -50233.156

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Hunk 3 (row group)
$ws.Range("H4").Value = 2988.25
$ws.Range("I4").Value = 101
$ws.Range("J4").Value = 3950.6667
$ws.Range("K4").Value = 101
$ws.Range("L4").Value = 3950.6667
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -4182.6667
# Hunk 4 (row group)
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("K5").Value = 150
$ws.Range("M5").Value = -38
# Hunk 5 (row group)
$ws.Range("H23").Value = 39266.668
# Hunk 6 (row group)
$ws.Range("H36").Value = 18194.666
$ws.Range("I36").Value = 7292
$ws.Range("K36").Value = 7292
$ws.Range("M36").Value = -6946
# Hunk 7 (row group)
$ws.Range("H37").Value = 29357.143
$ws.Range("I37").Value = 1500
$ws.Range("J37").Value = 34000
$ws.Range("K37").Value = 1500
$ws.Range("L37").Value = 34000
$ws.Range("M37").Value = -1227
$ws.Range("N37").Value = -34546
# Hunk 8 (row group)
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20976
# Hunk 9 (row group)
$ws.Range("H45").Value = 1466
$ws.Range("I45").Value = 865.6667
$ws.Range("J45").Value = 2666.6667
$ws.Range("K45").Value = 865.6667
$ws.Range("L45").Value = 2666.6667
$ws.Range("M45").Value = -488.6667
$ws.Range("N45").Value = -3420.6667
# Hunk 10 (row group)
$ws.Range("H61").Value = 52739948
$ws.Range("I61").Value = 71502680
$ws.Range("J61").Value = 204300
$ws.Range("K61").Value = 71502680
$ws.Range("L61").Value = 204300
$ws.Range("M61").Value = -71502468
$ws.Range("N61").Value = -204724
# Hunk 11 (row group)
$ws.Range("H80").Value = 34990
$ws.Range("I80").Value = 35000
$ws.Range("J80").Value = 34986.668
$ws.Range("K80").Value = 35000
$ws.Range("L80").Value = 34986.668
$ws.Range("M80").Value = -34002
$ws.Range("N80").Value = -36982.668
# Hunk 12 (row group)
$ws.Range("H83").Value = 34990
$ws.Range("I83").Value = 35000
$ws.Range("J83").Value = 34986.668
$ws.Range("K83").Value = 105000
$ws.Range("L83").Value = 104960.004
$ws.Range("M83").Value = -100008
$ws.Range("N83").Value = -114944.004
# Hunk 13 (row group)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Hunk 14 (row group)
$ws.Range("H132").Value = 36698.418
$ws.Range("I132").Value = 25993.707
$ws.Range("J132").Value = 59798.05
$ws.Range("K132").Value = 77981.121
$ws.Range("L132").Value = 179394.15
$ws.Range("M132").Value = -75451.121
$ws.Range("N132").Value = -184454.15
# Hunk 15 (row group)
$ws.Range("H136").Value = 52739948
$ws.Range("I136").Value = 71502680
$ws.Range("J136").Value = 204300
$ws.Range("K136").Value = 214508040
$ws.Range("L136").Value = 612900
$ws.Range("M136").Value = -214505490
$ws.Range("N136").Value = -618000

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Hunk 16 (row group)
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -35

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Hunk 17 (row group)
$ws.Range("H22").Value = 2976430.2
$ws.Range("I22").Value = 5208467
$ws.Range("J22").Value = 381
$ws.Range("K22").Value = 5208467
$ws.Range("L22").Value = 381
$ws.Range("M22").Value = -5208117
$ws.Range("N22").Value = -1081
# Hunk 18 (row group)
$ws.Range("H31").Value = 36708.1
$ws.Range("I31").Value = 22620.385
$ws.Range("J31").Value = 52633.348
$ws.Range("K31").Value = 22620.385
$ws.Range("L31").Value = 52633.348
$ws.Range("M31").Value = -22325.385
$ws.Range("N31").Value = -53223.348
# Hunk 19 (row group)
$ws.Range("H34").Value = 36708.1
$ws.Range("I34").Value = 22620.385
$ws.Range("J34").Value = 52633.348
$ws.Range("K34").Value = 22620.385
$ws.Range("L34").Value = 52633.348
$ws.Range("M34").Value = -22418.385
$ws.Range("N34").Value = -53037.348
# Hunk 20 (row group)
$ws.Range("H58").Value = 66668330
$ws.Range("I58").Value = 71430250
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 71430250
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -71430047
$ws.Range("N58").Value = -1906
# Hunk 21 (row group)
$ws.Range("H136").Value = 66668330
$ws.Range("I136").Value = 71430250
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 214290750
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -214288200
$ws.Range("N136").Value = -9600

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Hunk 22 (row group)
$ws.Range("H113").Value = 529.87177
$ws.Range("I113").Value = 482.0909
$ws.Range("J113").Value = 548.6429000000001
$ws.Range("K113").Value = 1446.2727
$ws.Range("L113").Value = 1645.9287
$ws.Range("M113").Value = 723.7273
$ws.Range("N113").Value = -5985.9287
# Hunk 23 (row group)
$ws.Range("H120").Value = 10227.143
$ws.Range("I120").Value = 5397.5
$ws.Range("J120").Value = 16666.666
$ws.Range("K120").Value = 16192.5
$ws.Range("L120").Value = 49999.99800000001
$ws.Range("M120").Value = -11354.5
$ws.Range("N120").Value = -59675.99800000001
# Hunk 24 (row group)
$ws.Range("H131").Value = 919.72974
$ws.Range("I131").Value = 399.16666
$ws.Range("J131").Value = 1020.4839
$ws.Range("K131").Value = 1197.49998
$ws.Range("L131").Value = 3061.4517
$ws.Range("M131").Value = 3842.50002
$ws.Range("N131").Value = -13141.4517

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Hunk 25 (row group)
$ws.Range("H2").Value = 128.38461
$ws.Range("I2").Value = 123.28571
$ws.Range("J2").Value = 134.33333
$ws.Range("K2").Value = 123.28571
$ws.Range("L2").Value = 134.33333
$ws.Range("M2").Value = -10.28570999999999
$ws.Range("N2").Value = -360.33333
# Hunk 26 (row group)
$ws.Range("H132").Value = 85769.46000000001
$ws.Range("I132").Value = 60897.65
$ws.Range("J132").Value = 146172.42
$ws.Range("K132").Value = 182692.95
$ws.Range("L132").Value = 438517.26
$ws.Range("M132").Value = -180162.95
$ws.Range("N132").Value = -443577.26

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Hunk 27 (row group)
$ws.Range("H82").Value = 1625.7142
$ws.Range("I82").Value = 1126.6666
$ws.Range("K82").Value = 1126.6666
$ws.Range("M82").Value = -765.6666
# Hunk 28 (row group)
$ws.Range("H85").Value = 1625.7142
$ws.Range("I85").Value = 1126.6666
$ws.Range("K85").Value = 1126.6666
$ws.Range("M85").Value = 121.3334
# Hunk 29 (row group)
$ws.Range("H93").Value = 1524.875
$ws.Range("I93").Value = 1524.875
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1524.875
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -276.875
$ws.Range("N93").ClearContents()
# Hunk 30 (row group)
$ws.Range("H132").Value = 65773.375
$ws.Range("I132").Value = 3216
$ws.Range("J132").Value = 203399.6
$ws.Range("K132").Value = 9648
$ws.Range("L132").Value = 610198.8
$ws.Range("M132").Value = -7118
$ws.Range("N132").Value = -615258.8
# Hunk 31 (row group)
$ws.Range("H136").Value = 34076.938
$ws.Range("I136").Value = 21955.75
$ws.Range("J136").Value = 86602.086
$ws.Range("K136").Value = 65867.25
$ws.Range("L136").Value = 259806.258
$ws.Range("M136").Value = -63317.25
$ws.Range("N136").Value = -264906.258

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Hunk 32 (row group)
$ws.Range("H133").Value = 39767
$ws.Range("J133").Value = 39767
$ws.Range("L133").Value = 39767
$ws.Range("N133").Value = -49887
# Hunk 33 (row group)
$ws.Range("H136").Value = 55707.633
$ws.Range("I136").Value = 44981.957
$ws.Range("J136").Value = 72153.664
$ws.Range("K136").Value = 134945.871
$ws.Range("L136").Value = 216460.992
$ws.Range("M136").Value = -132395.871
$ws.Range("N136").Value = -221560.992
